# Auto-generated edit script: update cryptos list values per diff
# commit: Updated cryptos list on Tue May 30 05:44:38 UTC 2023 with GitHub Actions

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to remain a text/string cell (matches the original
    # inlineStr cell type) even when the new value looks numeric, e.g.
    # "312.83" or "0.9080" - a plain .Value assignment would otherwise
    # let Excel auto-convert it to a real number and silently drop
    # formatting such as trailing zeros.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextValue $ws "D2" "27.820.92"
Set-TextValue $ws "E2" "  -1.05%  "
# Row 3
Set-TextValue $ws "D3" "1.903.20"
Set-TextValue $ws "E3" "  -0.36%  "
# Row 4
Set-TextValue $ws "E4" "  -0.11%  "
# Row 5
Set-TextValue $ws "D5" "312.83"
Set-TextValue $ws "E5" "  -0.98%  "
# Row 6
Set-TextValue $ws "E6" "  -0.06%  "
# Row 7
Set-TextValue $ws "D7" "0.5026"
Set-TextValue $ws "E7" "  +4.06%  "
# Row 8
Set-TextValue $ws "D8" "0.3812"
Set-TextValue $ws "E8" "  -0.19%  "
# Row 9
Set-TextValue $ws "D9" "0.07269"
Set-TextValue $ws "E9" "  -1.21%  "
# Row 10
Set-TextValue $ws "D10" "0.9080"
Set-TextValue $ws "E10" "  -2.92%  "
# Row 11
Set-TextValue $ws "D11" "20.81"
Set-TextValue $ws "E11" "  +0.07%  "
# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D12" "1.930.48"
Set-TextValue $ws "E12" "  +1.00%  "
# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D13" "0.07652"
Set-TextValue $ws "E13" "  -1.82%  "
# Row 14
Set-TextValue $ws "D14" "5.477"
Set-TextValue $ws "E14" "  -0.69%  "
# Row 15
Set-TextValue $ws "D15" "6.597"
Set-TextValue $ws "E15" "  -0.52%  "
# Row 16
Set-TextValue $ws "D16" "91.32"
Set-TextValue $ws "E16" "  -0.11%  "
# Row 17
Set-TextValue $ws "E17" "  -0.18%  "
# Row 18
Set-TextValue $ws "D18" "0.000008705"
Set-TextValue $ws "E18" "  -1.49%  "
# Row 19
Set-TextValue $ws "D19" "1.004"
Set-TextValue $ws "E19" "  +0.04%  "
# Row 20
Set-TextValue $ws "D20" "27.855.44"
# Row 21
Set-TextValue $ws "D21" "14.51"
Set-TextValue $ws "E21" "  -2.36%  "
# Row 22
Set-TextValue $ws "D22" "5.157"
Set-TextValue $ws "E22" "  +0.19%  "
# Row 23
Set-TextValue $ws "E23" "  -0.88%  "
# Row 24
Set-TextValue $ws "D24" "154.11"
Set-TextValue $ws "E24" "  -1.62%  "
# Row 25
Set-TextValue $ws "E25" "  -2.93%  "
# Row 26
Set-TextValue $ws "D26" "2.231"
Set-TextValue $ws "E26" "  +5.74%  "
# Row 27
Set-TextValue $ws "D27" "18.36"
Set-TextValue $ws "E27" "  -1.10%  "
# Row 28
Set-TextValue $ws "D28" "115.13"
Set-TextValue $ws "E28" "  -1.05%  "
# Row 29
Set-TextValue $ws "D29" "4.897"
Set-TextValue $ws "E29" "  -0.98%  "
# Row 30
Set-TextValue $ws "D30" "0.08972"
Set-TextValue $ws "E30" "  +0.58%  "
# Row 31
Set-TextValue $ws "D31" "3.204"
Set-TextValue $ws "E31" "  -4.28%  "
# Row 32
Set-TextValue $ws "D32" "1.230"
Set-TextValue $ws "E32" "  -1.85%  "
# Row 33
Set-TextValue $ws "D33" "0.7640"
# Row 34
Set-TextValue $ws "D34" "4.636"
Set-TextValue $ws "E34" "  -1.01%  "
# Row 35
Set-TextValue $ws "D35" "0.02053"
Set-TextValue $ws "E35" "  -0.09%  "
# Row 36
Set-TextValue $ws "D36" "2.539"
Set-TextValue $ws "E36" "  -2.88%  "
# Row 37
Set-TextValue $ws "D37" "1.096"
Set-TextValue $ws "E37" "  -0.67%  "
# Row 38
Set-TextValue $ws "D38" "0.5561"
Set-TextValue $ws "E38" "  +1.13%  "
# Row 39
Set-TextValue $ws "D39" "3.018"
Set-TextValue $ws "E39" "  +1.39%  "
# Row 40
Set-TextValue $ws "D40" "0.05251"
Set-TextValue $ws "E40" "  -1.24%  "
# Row 41
Set-TextValue $ws "D41" "6.969"
Set-TextValue $ws "E41" "  -0.78%  "
# Row 42
Set-TextValue $ws "D42" "8.471"
Set-TextValue $ws "E42" "  +0.20%  "
# Row 43
Set-TextValue $ws "D43" "0.1511"
Set-TextValue $ws "E43" "  -0.88%  "
# Row 44
Set-TextValue $ws "D44" "111.33"
Set-TextValue $ws "E44" "  +3.62%  "
# Row 45
Set-TextValue $ws "D45" "10.57"
Set-TextValue $ws "E45" "  -1.39%  "
# Row 46
Set-TextValue $ws "D46" "0.4787"
Set-TextValue $ws "E46" "  -0.96%  "
# Row 47
Set-TextValue $ws "E47" "  -0.05%  "
# Row 48
Set-TextValue $ws "E48" "  -1.74%  "
# Row 49
Set-TextValue $ws "D49" "67.34"
Set-TextValue $ws "E49" "  -1.60%  "
# Row 50
Set-TextValue $ws "D50" "0.06074"
Set-TextValue $ws "E50" "  -0.58%  "
# Row 51
Set-TextValue $ws "E51" "  -0.44%  "
